$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Annotation values (column E, "label") for rows 34-54, continuing the
# human annotation that previously only went up through row 33.
$values = @(-1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, -2, 0, 0, -1, 0)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 34 + $i
    $ws.Range("E$row").Value = $values[$i]
}

# Move the selection/active cell to E55, matching the cursor position left
# behind after annotating through row 54.
$ws.Range("A54").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 54
$ws.Range("E55").Select() | Out-Null
